$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# 1. Update the "last updated" timestamp text (shared string change)
$ws.Range("A1").Value = "Datos actualizados a 28 de Julio de 2020 a las 06:18"

# 2. Brasil row (row 5): update case counters
$ws.Range("B5").Value = 2446397
$ws.Range("C5").Value = 2917
$ws.Range("E5").Value = 690993
$ws.Range("G5").Value = 58
$ws.Range("H5").Value = 87737

# 3. Guyana overtakes Bahamas and Burundi in the ranking (rows 164-166 swap
#    countries while keeping the row positions sorted by total cases).
#    Row 164 becomes Guyana, row 165 becomes Bahamas, row 166 becomes Burundi.
$ws.Range("A164").Value = "Guyana"
$ws.Range("B164").Value = 389
$ws.Range("C164").Value = 0
$ws.Range("D164").Value = 181
$ws.Range("E164").Value = 188
$ws.Range("F164").Value = 0
$ws.Range("G164").Value = 0
$ws.Range("H164").Value = 20

$ws.Range("A165").Value = "Bahamas"
$ws.Range("B165").Value = 382
$ws.Range("C165").Value = 0
$ws.Range("D165").Value = 91
$ws.Range("E165").Value = 280
$ws.Range("F165").Value = 0
$ws.Range("G165").Value = 0
$ws.Range("H165").Value = 11

$ws.Range("A166").Value = "Burundi"
$ws.Range("B166").Value = 378
$ws.Range("C166").Value = 0
$ws.Range("D166").Value = 301
$ws.Range("E166").Value = 76
$ws.Range("F166").Value = 0
$ws.Range("G166").Value = 0
$ws.Range("H166").Value = 1

# 4. Mongolia (row 172): update case counters
$ws.Range("B172").Value = 289
$ws.Range("C172").Value = 1
$ws.Range("D172").Value = 222
$ws.Range("E172").Value = 67

# 5. Islas Turcas y Caicos (row 188): update case counters
$ws.Range("B188").Value = 99
$ws.Range("D188").Value = 36
$ws.Range("E188").Value = 61
